# Updates the "Price" column (column D) of the crypto symbol list with the
# latest scraped values, as produced by the scheduled GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new price text. NumberFormat is forced to "@" (Text)
# before assigning so the value is stored as text (matching the sheet's
# existing inline-string / text cells) instead of being auto-coerced into a
# floating point number by Excel's "General" numeric auto-detection.
$updates = [ordered]@{
    "D2"  = "282.58"
    "D3"  = "20.61"
    "D5"  = "0.06174"
    "D6"  = "3.582"
    "D7"  = "6.559"
    "D8"  = "1.499"
    "D9"  = "0.8186"
    "D11" = "0.1647"
    "D12" = "0.08425"
    "D13" = "0.03523"
    "D14" = "0.03216"
    "D15" = "0.09146"
    "D16" = "3.721"
    "D17" = "0.001638"
    "D18" = "0.04706"
    "D19" = "0.006562"
    "D20" = "0.006171"
    "D22" = "0.0001602"
    "D23" = "3.835"
    "D41" = "0.007196"
    "D42" = "0.1098"
    "D43" = "0.004505"
    "D44" = "0.01142"
    "D45" = "0.00006609"
    "D47" = "1.052"
    "D48" = "0.002848"
    "D49" = "0.00001902"
    "D50" = "0.01241"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
